$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, [string]$value) {
    # Force the cell to store the value as literal text even when it
    # "looks" numeric (e.g. "123"), without leaving a lingering custom
    # number-format style behind on the cell.
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

# Row 2 - Luiz Antônio de Souza
$ws.Range("A2").Value = "Luiz Antônio de Souza"
$ws.Range("B2").Value = "408.391.907-32"
$ws.Range("C2").Value = "luiz.santos98@gmail.com"
$ws.Range("D2").Value = "(11) 95234-6781"
$ws.Range("E2").Value = "Rua das Flores, 123, Bairro Jardim, São Paulo - SP, CEP: 01234-567"

# Row 3 - Luis Gustavo
$ws.Range("A3").Value = "Luis Gustavo"
$ws.Range("B3").Value = "782.594.315-06"
$ws.Range("C3").Value = "gustavo.martins86@yahoo.com.br"
$ws.Range("D3").Value = "(31) 98712-3456"
$ws.Range("E3").Value = "Rua Afonso Pena, 789, Centro, Belo Horizonte - MG, CEP: 30130-907"

# Row 4 - Renan Carlos
$ws.Range("A4").Value = "Renan Carlos"
$ws.Range("B4").Value = "145.782.963-14"
$ws.Range("C4").Value = "renan.alves99@outlook.com"
$ws.Range("D4").Value = "(41) 99876-5432"
$ws.Range("E4").Value = "Rua das Palmeiras, 100, Batel, Curitiba - PR, CEP: 80240-001"

# Row 5 - João Davi
$ws.Range("A5").Value = "João Davi"
$ws.Range("B5").Value = "613.809.124-50"
$ws.Range("C5").Value = "joao.silva23@hotmail.com"
$ws.Range("D5").Value = "(21) 98765-4321"
$ws.Range("E5").Value = "Avenida Atlântica, 456, Copacabana, Rio de Janeiro - RJ, CEP: 22070-001"

# Row 6 (new) - fasfas
$ws.Range("A6").Value = "fasfas"
Set-TextValue $ws.Range("B6") "123"
$ws.Range("C6").Value = "frwef"
Set-TextValue $ws.Range("D6") "23131"
$ws.Range("E6").Value = "fsdf"
